$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change Y27's style from date-only to date-time format, keep the same value
$ws.Range("Y27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("Y27").Value = 45757

# Add new row 28 with values
$rowValues = @{
    "A28" = 565
    "B28" = 479
    "C28" = 450
    "D28" = 548
    "E28" = 498
    "F28" = 532
    "G28" = 474
    "H28" = 567
    "I28" = 495
    "J28" = 450
    "K28" = 572
    "L28" = 483
    "M28" = 462
    "N28" = 505
    "O28" = 557
    "P28" = 483
    "Q28" = 618
    "R28" = 497
    "S28" = 474
    "T28" = 485
    "U28" = 619
    "V28" = 540
    "W28" = 594
    "X28" = 490
    "Y28" = 45758
    "Z28" = 818
    "AA28" = 556
    "AB28" = 525.5
    "AC28" = 512
    "AD28" = 542
    "AE28" = 504
    "AF28" = 506
    "AG28" = 730
    "AH28" = 467
    "AI28" = 720
    "AJ28" = 474
    "AK28" = 484
    "AL28" = 550
    "AM28" = 540
    "AN28" = 485
    "AO28" = 540
    "AP28" = 529
    "AQ28" = 565
    "AR28" = 546
    "AS28" = 635
    "AT28" = 637
    "AU28" = 493
    "AV28" = 475
}

foreach ($addr in $rowValues.Keys) {
    $ws.Range($addr).Value = $rowValues[$addr]
}

# Y28 should be formatted as date-only (s=3), matching Y27's original style
$ws.Range("Y28").NumberFormat = "YYYY-MM-DD"
